$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "Datos actualizados..." timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 21:50"

# Update row 29 (Leon) statistics
$ws.Range("B29").Value = 1145
$ws.Range("C29").Value = 90
$ws.Range("D29").Value = 1009
$ws.Range("E29").Value = 46
